$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.106.90'
$ws.Range("E2").Value = '  +0.94%  '
$ws.Range("D3").Value = '1.636.98'
$ws.Range("E3").Value = '  -0.09%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '216.77'
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("D6").Value = '0.517'
$ws.Range("E6").Value = '  +1.76%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '0.254'
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").Value = '0.0625'
$ws.Range("E9").Value = '  +0.22%  '
$ws.Range("D10").Value = '19.92'
$ws.Range("E10").Value = '  +0.30%  '
$ws.Range("D11").Value = '0.0847'
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("D12").Value = '1.866.40'
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("D13").Value = '1.641.19'
$ws.Range("E13").Value = '  +0.10%  '
$ws.Range("E14").Value = '  +0.16%  '
$ws.Range("D15").Value = '0.541'
$ws.Range("E15").Value = '  +2.17%  '
$ws.Range("D16").Value = '66.68'
$ws.Range("E16").Value = '  -0.63%  '
$ws.Range("D17").Value = '27.092.60'
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("E18").Value = '  +1.16%  '
$ws.Range("D19").Value = '216.90'
$ws.Range("E19").Value = '  -0.62%  '
$ws.Range("E20").Value = '  +0.05%  '
$ws.Range("E21").Value = '  +2.05%  '
$ws.Range("E22").Value = '  +0.28%  '
$ws.Range("D23").Value = '2.49'
$ws.Range("E23").Value = '  +1.96%  '
$ws.Range("E24").Value = '  -0.75%  '
$ws.Range("D25").Value = '146.68'
$ws.Range("E25").Value = '  -0.38%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("D27").Value = '7.40'
$ws.Range("E27").Value = '  +1.97%  '
$ws.Range("E28").Value = '  +0.24%  '
$ws.Range("D29").Value = '15.67'
$ws.Range("E29").Value = '  -0.46%  '
$ws.Range("E30").Value = '  +0.87%  '
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("E32").Value = '  +1.40%  '
$ws.Range("E33").Value = '  +0.51%  '
$ws.Range("D34").Value = '1.300.62'
$ws.Range("E34").Value = '  +2.78%  '
$ws.Range("E35").Value = '  +0.16%  '
$ws.Range("E37").Value = '  -0.51%  '
$ws.Range("D38").Value = '0.855'
$ws.Range("E38").Value = '  +2.34%  '
$ws.Range("D39").Value = '0.542'
$ws.Range("E39").Value = '  +1.67%  '
$ws.Range("E40").Value = '  +0.07%  '
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("E42").Value = '  +5.64%  '
$ws.Range("E43").Value = '  -1.39%  '
$ws.Range("D44").Value = '1.776.53'
$ws.Range("E44").Value = '  -0.17%  '
$ws.Range("D45").Value = '61.68'
$ws.Range("E45").Value = '  -0.48%  '
$ws.Range("D46").Value = '91.15'
$ws.Range("E46").Value = '  -0.72%  '
$ws.Range("E47").Value = '  +0.93%  '
$ws.Range("D48").Value = '0.0₆0108'
$ws.Range("E48").Value = '  +2.09%  '
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("D50").Value = '7.64'
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("D51").Value = '0.0957'
$ws.Range("E51").Value = '  -0.25%  '
